$wb = $excel.ActiveWorkbook

# --- 1. Insert a new "2022-Q1" sheet right before the "总计" (Total) sheet ---
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# NOTE: $totalSheet was resolved to a position-based anchor; once the new
# sheet is spliced in *before* it, that same anchor now resolves to the
# freshly inserted sheet instead of following "总计" to its new position.
# Re-resolve it by name so later edits land on the right sheet.
$totalSheet = $wb.Worksheets.Item("总计")

# Copy layout/formatting from the "2021-Q4" sheet (same column layout: A..H)
$srcSheet = $wb.Worksheets.Item("2021-Q4")
$srcSheet.Range("A1:H4").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)
$newSheet.Range("A1").Clear()

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row 2
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'217024"
$newSheet.Range("C2").Value = "招商安盈债券"
$newSheet.Range("D2").Value = "'35.05"
$newSheet.Range("E2").Value = "'20.20"
$newSheet.Range("F2").Value = "'3.72"
$newSheet.Range("G2").Value = "'1.3039"
$newSheet.Range("H2").Value = 1

# Row 3
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'159611"
$newSheet.Range("C3").Value = "广发中证全指电力ETF"
$newSheet.Range("D3").Value = "'13.38"
$newSheet.Range("E3").Value = "'99.14"
$newSheet.Range("F3").Value = "'3.47"
$newSheet.Range("G3").Value = "'0.4643"
$newSheet.Range("H3").Value = 7

# Row 4
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'014887"
$newSheet.Range("C4").Value = "招商安福1年定期开放债券"
$newSheet.Range("D4").Value = "'17.22"
$newSheet.Range("E4").Value = "'27.65"
$newSheet.Range("F4").Value = "'2.40"
$newSheet.Range("G4").Value = "'0.4133"
$newSheet.Range("H4").Value = 2

# The leading apostrophes above force text-storage for numeric-looking strings
# (preserving e.g. the "014887" leading zero); ClearFormats drops the
# resulting quote-prefix style so the cells stay visually identical to their
# un-styled neighbours on the other sheets.
$newSheet.Range("B2:G4").ClearFormats()

# --- 2. Prepend a "2022-Q1" row to the "总计" (Total) summary sheet ---
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 2.18

# Renumber the index column (A) for the rows that shifted down
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5

Write-Host "done"
